$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4999.8335
$ws.Range("I40").Value = 4999
$ws.Range("J40").Value = 4999.909
$ws.Range("K40").Value = 4999
$ws.Range("L40").Value = 4999.909
$ws.Range("M40").Value = -4824
$ws.Range("N40").Value = -5349.909
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H55").Value = 402.05554
$ws.Range("J55").Value = 433.5625
$ws.Range("L55").Value = 433.5625
$ws.Range("N55").Value = -861.5625
$ws.Range("H100").Value = 9102.762000000001
$ws.Range("I100").Value = 899.8333
$ws.Range("K100").Value = 899.8333
$ws.Range("M100").Value = -358.8333
$ws.Range("H105").Value = 84999.5
$ws.Range("J105").Value = 84999.5
$ws.Range("L105").Value = 84999.5
$ws.Range("N105").Value = -91987.5
$ws.Range("H113").Value = 7916.6665
$ws.Range("I113").Value = 5833.3335
$ws.Range("K113").Value = 5833.3335
$ws.Range("M113").Value = -2579.3335
$ws.Range("H132").Value = 50016770
$ws.Range("I132").Value = 83343450
$ws.Range("J132").Value = 26742.75
$ws.Range("K132").Value = 250030350
$ws.Range("L132").Value = 80228.25
$ws.Range("M132").Value = -250027820
$ws.Range("N132").Value = -85288.25
$ws.Range("H138").Value = 324966.78
$ws.Range("I138").Value = 37939.816
$ws.Range("K138").Value = 113819.448
$ws.Range("M138").Value = -108679.448
$ws.Range("H139").Value = 84000
$ws.Range("J139").Value = 84000
$ws.Range("L139").Value = 84000
$ws.Range("N139").Value = -94280

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 12831.75
$ws.Range("I5").Value = 12831.75
$ws.Range("K5").Value = 12831.75
$ws.Range("M5").Value = -12719.75
$ws.Range("H11").Value = 15996.75
$ws.Range("I11").Value = 9999
$ws.Range("J11").Value = 21994.5
$ws.Range("K11").Value = 9999
$ws.Range("L11").Value = 21994.5
$ws.Range("M11").Value = -9855
$ws.Range("N11").Value = -22282.5
$ws.Range("H32").Value = 5590.5156
$ws.Range("I32").Value = 5304.933
$ws.Range("K32").Value = 5304.933
$ws.Range("M32").Value = -5017.933
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H124").Value = 87249.75
$ws.Range("J124").Value = 87249.75
$ws.Range("L124").Value = 87249.75
$ws.Range("N124").Value = -97069.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 12831.75
$ws.Range("I4").Value = 12831.75
$ws.Range("K4").Value = 12831.75
$ws.Range("M4").Value = -12716.75
$ws.Range("H22").Value = 199.7
$ws.Range("I22").Value = 92.166664
$ws.Range("J22").Value = 361
$ws.Range("K22").Value = 92.166664
$ws.Range("L22").Value = 361
$ws.Range("M22").Value = 80.833336
$ws.Range("N22").Value = -707
$ws.Range("H94").Value = 2228.9033
$ws.Range("J94").Value = 3142.125
$ws.Range("L94").Value = 3142.125
$ws.Range("N94").Value = -4044.125
$ws.Range("H134").Value = 5621.576
$ws.Range("I134").Value = 2014.4828
$ws.Range("K134").Value = 6043.4484
$ws.Range("M134").Value = -3508.4484

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 51925
$ws.Range("I59").Value = 10000
$ws.Range("K59").Value = 10000
$ws.Range("M59").Value = -8855
$ws.Range("H62").Value = 3671.2222
$ws.Range("I62").Value = 3637.25
$ws.Range("J62").Value = 3698.4
$ws.Range("K62").Value = 3637.25
$ws.Range("L62").Value = 3698.4
$ws.Range("M62").Value = -3013.25
$ws.Range("N62").Value = -4946.4
$ws.Range("H65").Value = 3671.2222
$ws.Range("I65").Value = 3637.25
$ws.Range("J65").Value = 3698.4
$ws.Range("K65").Value = 18186.25
$ws.Range("L65").Value = 18492
$ws.Range("M65").Value = -15066.25
$ws.Range("N65").Value = -24732
$ws.Range("H94").Value = 1222.7
$ws.Range("I94").Value = 1129.6666
$ws.Range("J94").Value = 1239.1177
$ws.Range("K94").Value = 1129.6666
$ws.Range("L94").Value = 1239.1177
$ws.Range("M94").Value = -678.6666
$ws.Range("N94").Value = -2141.1177
$ws.Range("H99").Value = 2597
$ws.Range("I99").Value = 3697.5
$ws.Range("J99").Value = 1496.5
$ws.Range("K99").Value = 3697.5
$ws.Range("L99").Value = 1496.5
$ws.Range("M99").Value = -2199.5
$ws.Range("N99").Value = -4492.5
$ws.Range("H126").Value = 2597
$ws.Range("I126").Value = 3697.5
$ws.Range("J126").Value = 1496.5
$ws.Range("K126").Value = 11092.5
$ws.Range("L126").Value = 4489.5
$ws.Range("M126").Value = -8622.5
$ws.Range("N126").Value = -9429.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 683.5833
$ws.Range("I86").Value = 682.0909
$ws.Range("J86").Value = 700
$ws.Range("K86").Value = 2046.2727
$ws.Range("L86").Value = 2100
$ws.Range("M86").Value = -860.2727
$ws.Range("N86").Value = -4472
$ws.Range("H89").Value = 683.5833
$ws.Range("I89").Value = 682.0909
$ws.Range("J89").Value = 700
$ws.Range("K89").Value = 6138.8181
$ws.Range("L89").Value = 6300
$ws.Range("M89").Value = -210.8181000000004
$ws.Range("N89").Value = -18156
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H121").Value = 3099.75
$ws.Range("I121").Value = 744.25
$ws.Range("J121").Value = 7810.75
$ws.Range("K121").Value = 2232.75
$ws.Range("L121").Value = 23432.25
$ws.Range("M121").Value = -922.75
$ws.Range("N121").Value = -26052.25
$ws.Range("H137").Value = 1940.3478
$ws.Range("J137").Value = 2932.1667
$ws.Range("L137").Value = 8796.500100000001
$ws.Range("N137").Value = -18996.5001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 27777
$ws.Range("J24").Value = 27777
$ws.Range("L24").Value = 27777
$ws.Range("N24").Value = -28123
$ws.Range("H34").Value = 162941.4
$ws.Range("J34").Value = 173516.67
$ws.Range("L34").Value = 173516.67
$ws.Range("N34").Value = -174052.67
$ws.Range("H76").Value = 162941.4
$ws.Range("J76").Value = 173516.67
$ws.Range("L76").Value = 173516.67
$ws.Range("N76").Value = -174146.67
$ws.Range("H79").Value = 162941.4
$ws.Range("J79").Value = 173516.67
$ws.Range("L79").Value = 173516.67
$ws.Range("N79").Value = -175700.67

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2137.25
$ws.Range("I40").Value = 2149.7727
$ws.Range("K40").Value = 2149.7727
$ws.Range("M40").Value = -2013.7727
$ws.Range("H46").Value = 4455.35
$ws.Range("J46").Value = 5702
$ws.Range("L46").Value = 5702
$ws.Range("N46").Value = -6078
$ws.Range("H55").Value = 197.10527
$ws.Range("I55").Value = 99.75
$ws.Range("J55").Value = 267.9091
$ws.Range("K55").Value = 99.75
$ws.Range("L55").Value = 267.9091
$ws.Range("M55").Value = 73.25
$ws.Range("N55").Value = -613.9091000000001
$ws.Range("H68").Value = 2499.8
$ws.Range("J68").Value = 2499.6667
$ws.Range("L68").Value = 2499.6667
$ws.Range("N68").Value = -3997.6667
$ws.Range("H71").Value = 2499.8
$ws.Range("J71").Value = 2499.6667
$ws.Range("L71").Value = 12498.3335
$ws.Range("N71").Value = -19986.3335
$ws.Range("H120").Value = 103110
$ws.Range("J120").Value = 103110
$ws.Range("L120").Value = 103110
$ws.Range("N120").Value = -112786
$ws.Range("H122").Value = 4098.364
$ws.Range("I122").Value = 3926.6667
$ws.Range("J122").Value = 4217.231
$ws.Range("K122").Value = 11780.0001
$ws.Range("L122").Value = 12651.693
$ws.Range("M122").Value = -9330.000100000001
$ws.Range("N122").Value = -17551.693

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 45999.5
$ws.Range("J69").Value = 45999.5
$ws.Range("L69").Value = 45999.5
$ws.Range("N69").Value = -47497.5
$ws.Range("H72").Value = 45999.5
$ws.Range("J72").Value = 45999.5
$ws.Range("L72").Value = 137998.5
$ws.Range("N72").Value = -145486.5
